$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row (row 3) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Retour status"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("D3").Value = "Klantenservice / Opvolging"
$logs.Range("F3").Value = "2025-08-26 19:33:14"
$logs.Range("G3").Value = "Nee"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# --- Extend the conditional formatting ranges to cover the new row ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $fcs = $logs.Range("$col`2").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range("$col`2:$col`3"))
    }
}

# --- Dashboard sheet: bump the count for "Klantenservice / Opvolging" ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 2
